$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7435.946
$ws.Range("I32").Value = 1715.091
$ws.Range("K32").Value = 1715.091
$ws.Range("M32").Value = -1428.091

$ws.Range("H45").Value = 3749.9412
$ws.Range("I45").Value = 2583
$ws.Range("K45").Value = 2583
$ws.Range("M45").Value = -2206

$ws.Range("H74").Value = 5068.857
$ws.Range("I74").Value = 3997
$ws.Range("K74").Value = 3997
$ws.Range("M74").Value = -3123

$ws.Range("H77").Value = 5068.857
$ws.Range("I77").Value = 3997
$ws.Range("K77").Value = 19985
$ws.Range("M77").Value = -15617

$ws.Range("H110").Value = 2390.6155
$ws.Range("I110").Value = 2120
$ws.Range("K110").Value = 2120
$ws.Range("M110").Value = -75

$ws.Range("H132").Value = 2948.7334
$ws.Range("I132").Value = 2659.3572
$ws.Range("K132").Value = 7978.071599999999
$ws.Range("M132").Value = -5448.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 9999.5
$ws.Range("I33").Value = 9999.5
$ws.Range("K33").Value = 9999.5
$ws.Range("M33").Value = -9663.5

$ws.Range("H103").Value = 18999.5
$ws.Range("J103").Value = 18999.5
$ws.Range("L103").Value = 18999.5
$ws.Range("N103").Value = -21343.5

$ws.Range("H134").Value = 3984.2856
$ws.Range("I134").Value = 3998.3333
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 11994.9999
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -9459.999899999999
$ws.Range("N134").Value = -16770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1006.8
$ws.Range("I22").Value = 989.7857
$ws.Range("K22").Value = 989.7857
$ws.Range("M22").Value = -639.7857

$ws.Range("H35").Value = 1095.5555
$ws.Range("I35").Value = 1095.5555
$ws.Range("K35").Value = 1095.5555
$ws.Range("M35").Value = -801.5554999999999

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H134").Value = 1193.8422
$ws.Range("I134").Value = 1093.5
$ws.Range("J134").Value = 1365.8572
$ws.Range("K134").Value = 3280.5
$ws.Range("L134").Value = 4097.571599999999
$ws.Range("M134").Value = -745.5
$ws.Range("N134").Value = -9167.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1316.8
$ws.Range("J68").Value = 1396
$ws.Range("L68").Value = 4188
$ws.Range("N68").Value = -5810

$ws.Range("H71").Value = 1316.8
$ws.Range("J71").Value = 1396
$ws.Range("L71").Value = 12564
$ws.Range("N71").Value = -20676

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H107").Value = 459.70587
$ws.Range("I107").Value = 525.5
$ws.Range("J107").Value = 423.81818
$ws.Range("K107").Value = 1576.5
$ws.Range("L107").Value = 1271.45454
$ws.Range("M107").Value = 343.5
$ws.Range("N107").Value = -5111.45454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 22361
$ws.Range("J96").Value = 22361
$ws.Range("L96").Value = 22361
$ws.Range("N96").Value = -27853

$ws.Range("H126").Value = 7907.3477
$ws.Range("I126").Value = 6964.3887
$ws.Range("J126").Value = 11302
$ws.Range("K126").Value = 20893.1661
$ws.Range("L126").Value = 33906
$ws.Range("M126").Value = -18423.1661
$ws.Range("N126").Value = -38846

$ws.Range("H132").Value = 5978.516
$ws.Range("I132").Value = 5666.3794
$ws.Range("K132").Value = 16999.1382
$ws.Range("M132").Value = -14469.1382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H41").Value = 22999
$ws.Range("I41").Value = 22999
$ws.Range("K41").Value = 22999
$ws.Range("M41").Value = -22561

$ws.Range("H48").Value = 34721.5
$ws.Range("I48").Value = 29629
$ws.Range("K48").Value = 29629
$ws.Range("M48").Value = -28968

$ws.Range("H93").Value = 1580.56
$ws.Range("I93").Value = 1512.2941
$ws.Range("K93").Value = 1512.2941
$ws.Range("M93").Value = -264.2941000000001

$ws.Range("H132").Value = 13143
$ws.Range("I132").Value = 13143
$ws.Range("K132").Value = 39429
$ws.Range("M132").Value = -36899

$ws.Range("H136").Value = 6276.5454
$ws.Range("I136").Value = 7377.4287
$ws.Range("J136").Value = 4350
$ws.Range("K136").Value = 22132.2861
$ws.Range("L136").Value = 13050
$ws.Range("M136").Value = -19582.2861
$ws.Range("N136").Value = -18150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3516.6667
$ws.Range("I17").Value = 3550
$ws.Range("K17").Value = 3550
$ws.Range("M17").Value = -3378

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H55").Value = 10048
$ws.Range("I55").Value = 10048
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 10048
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -9771
$ws.Range("N55").ClearContents()

$ws.Range("H126").Value = 3122.375
$ws.Range("I126").Value = 2139.8572
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 6419.571599999999
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -3949.571599999999
$ws.Range("N126").Value = -34940

$ws.Range("H136").Value = 15597.2
$ws.Range("I136").Value = 20995.334
$ws.Range("K136").Value = 62986.00199999999
$ws.Range("M136").Value = -60436.00199999999

$ws.Range("H137").Value = 79390.2
$ws.Range("J137").Value = 79390.2
$ws.Range("L137").Value = 79390.2
$ws.Range("N137").Value = -89590.2
